$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# "Ready for handoff" generate/handoff timestamp moved forward (18:32:38 -> 18:33:13).
# Referenced by Overview!G4:G7 and de-de!H4:H7.
foreach ($row in 4..7) {
    $wsOverview.Cells.Item($row, 7).Value = "2016-08-22 18:33:13"
    $wsDe.Cells.Item($row, 8).Value = "2016-08-22 18:33:13"
}

# zh-cn Latest Handoff Datetime moved forward (18:32:31 -> 18:33:03), referenced by zh-cn!H4:H7.
foreach ($row in 4..7) {
    $wsZh.Cells.Item($row, 8).Value = "2016-08-22 18:33:03"
}

# Priority recalculated from "low" to "ht" for rows 4-7 in both language sheets.
foreach ($row in 4..7) {
    $wsZh.Cells.Item($row, 5).Value = "ht"
    $wsDe.Cells.Item($row, 5).Value = "ht"
}
